$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# TC2's step #2 (row 20) and TC3's step #2 (row 28) need to swap their
# "Steps" (col B) and "Expected Results" (col D) content, so that the
# "realizar a liquidacao" step now belongs to TC2 and the
# "atribuir/desatribuir" step now belongs to TC3.

$tc2StepB = $ws.Range("B20").Value2
$tc2StepD = $ws.Range("D20").Value2
$tc3StepB = $ws.Range("B28").Value2
$tc3StepD = $ws.Range("D28").Value2

$ws.Range("B20").Value = $tc3StepB
$ws.Range("D20").Value = $tc3StepD

$ws.Range("B28").Value = $tc2StepB
$ws.Range("D28").Value = $tc2StepD
